$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells appended after the existing J1 ("longitud") header.
$ws.Range("K1").Value = "fecha"
$ws.Range("L1").Value = "conductor "
$ws.Range("M1").Value = "zona"
$ws.Range("N1").Value = "n_entregas"
$ws.Range("O1").Value = "tiempo_total"
$ws.Range("P1").Value = "combustible_usado"
$ws.Range("Q1").Value = "km_recorridos"

# Match the header formatting already used by A1:J1 (style index 1).
$ws.Range("A1").Copy()
$ws.Range("K1:Q1").PasteSpecial(-4122)

# Widen the new "combustible_usado" column like the others.
$ws.Columns.Item(16).ColumnWidth = 22
